# Updated cryptos list on Thu Dec 21 11:47:39 UTC 2023 with GitHub Actions
#
# Refresh the crypto price/volume table on the active sheet with the
# latest scraped values (including two rows whose ranking swapped with
# a neighbor). Numeric-looking text values (e.g. "1.00", "0.0930") are
# written through a brief NumberFormat "@" (Text) toggle so Excel does
# not reinterpret them as numbers -- the format is restored to
# "General" immediately after, matching the original cell formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.965.72"
$ws.Range("E2").Value = "  +2.58%  "
$ws.Range("D3").Value = "2.251.52"
$ws.Range("E3").Value = "  +1.91%  "
$ws.Range("E4").Value = "  +0.09%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "270.39"
$c.NumberFormat = "General"
$ws.Range("E5").Value = "  +5.73%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "87.64"
$c.NumberFormat = "General"
$ws.Range("E6").Value = "  +13.53%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.624"
$c.NumberFormat = "General"
$ws.Range("E7").Value = "  +2.10%  "
$ws.Range("E8").Value = "  +0.02%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.614"
$c.NumberFormat = "General"
$ws.Range("E9").Value = "  +3.85%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "45.82"
$c.NumberFormat = "General"
$ws.Range("E10").Value = "  +7.37%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.0930"
$c.NumberFormat = "General"
$ws.Range("E11").Value = "  +2.89%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "7.58"
$c.NumberFormat = "General"
$ws.Range("E12").Value = "  +8.53%  "
$ws.Range("E13").Value = "  +2.83%  "
$ws.Range("D14").Value = "2.583.68"
$ws.Range("E14").Value = "  +1.72%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "15.04"
$c.NumberFormat = "General"
$ws.Range("E15").Value = "  +4.71%  "
$ws.Range("D16").Value = "2.262.70"
$ws.Range("E16").Value = "  +2.97%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.799"
$c.NumberFormat = "General"
$ws.Range("E17").Value = "  +2.13%  "
$ws.Range("D18").Value = "43.952.85"
$ws.Range("E18").Value = "  +2.71%  "
$ws.Range("E19").Value = "  +0.50%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "6.05"
$c.NumberFormat = "General"
$ws.Range("E20").Value = "  +1.63%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "70.30"
$c.NumberFormat = "General"
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "2.40"
$c.NumberFormat = "General"
$ws.Range("E22").Value = "  +3.56%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "233.82"
$c.NumberFormat = "General"
$ws.Range("E23").Value = "  +1.73%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "8.80"
$c.NumberFormat = "General"
$ws.Range("E24").Value = "  -4.38%  "
$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.57"
$c.NumberFormat = "General"
$ws.Range("E25").Value = "  +17.16%  "
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.NumberFormat = "General"
$ws.Range("E26").Value = "  +0.00%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "11.00"
$c.NumberFormat = "General"
$ws.Range("E27").Value = "  +2.87%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "3.57"
$c.NumberFormat = "General"
$ws.Range("E28").Value = "  +6.76%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "40.66"
$c.NumberFormat = "General"
$ws.Range("E29").Value = "  -4.51%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "2.34"
$c.NumberFormat = "General"
$ws.Range("E30").Value = "  +5.99%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "175.20"
$c.NumberFormat = "General"
$ws.Range("E31").Value = "  +1.66%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.0913"
$c.NumberFormat = "General"
$ws.Range("E32").Value = "  +3.83%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "20.79"
$c.NumberFormat = "General"
$ws.Range("E33").Value = "  +2.26%  "
$ws.Range("E34").Value = "  +4.24%  "
$ws.Range("E35").Value = "  +2.21%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.113"
$c.NumberFormat = "General"
$ws.Range("E36").Value = "  +6.25%  "
$ws.Range("E37").Value = "  -0.23%  "
$ws.Range("E38").Value = "  +1.00%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "3.38"
$c.NumberFormat = "General"
$ws.Range("E39").Value = "  +16.99%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "12.84"
$c.NumberFormat = "General"
$ws.Range("E40").Value = "  -1.91%  "
$ws.Range("E41").Value = "  +2.37%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "65.33"
$c.NumberFormat = "General"
$ws.Range("E42").Value = "  +7.14%  "
$ws.Range("E43").Value = "  +2.57%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "5.42"
$c.NumberFormat = "General"
$ws.Range("E44").Value = "  +2.50%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.0991"
$c.NumberFormat = "General"
$ws.Range("E45").Value = "  +2.06%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "8.41"
$c.NumberFormat = "General"
$ws.Range("E46").Value = "  -0.42%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "100.55"
$c.NumberFormat = "General"
$ws.Range("E47").Value = "  -2.29%  "
$ws.Range("E48").Value = "  +7.71%  "
$ws.Range("E49").Value = "  +2.61%  "
$ws.Range("B50").Value = "WOONetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.440"
$c.NumberFormat = "General"
$ws.Range("E50").Value = "  -4.96%  "
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "1.51"
$c.NumberFormat = "General"
$ws.Range("E51").Value = "  +1.10%  "
